{"js": "// Replace the date and each \"a\u00f7b=\" exercise text with its updated value.\n// The mapping below mirrors the unified diff 1:1, in document order, and\n// every \"old\" value is unique within the document, so a simple\n// search-and-replace per pair is unambiguous and safe.\nconst replacements = [\n  [\"2024-05-02 Thursday\", \"2024-05-03 Friday\"],\n  [\"261\u00f73=\", \"250\u00f75=\"],\n  [\"618\u00f76=\", \"977\u00f74=\"],\n  [\"811\u00f77=\", \"931\u00f78=\"],\n  [\"436\u00f76=\", \"446\u00f77=\"],\n  [\"796\u00f76=\", \"380\u00f75=\"],\n  [\"746\u00f79=\", \"818\u00f79=\"],\n  [\"113\u00f74=\", \"715\u00f75=\"],\n  [\"351\u00f78=\", \"332\u00f77=\"],\n  [\"284\u00f74=\", \"496\u00f76=\"],\n  [\"904\u00f79=\", \"492\u00f77=\"],\n  [\"328\u00f78=\", \"293\u00f74=\"],\n  [\"199\u00f75=\", \"960\u00f72=\"],\n  [\"749\u00f79=\", \"749\u00f75=\"],\n  [\"301\u00f73=\", \"192\u00f79=\"],\n  [\"160\u00f72=\", \"133\u00f75=\"],\n  [\"851\u00f77=\", \"181\u00f77=\"],\n  [\"511\u00f77=\", \"571\u00f73=\"],\n  [\"276\u00f75=\", \"165\u00f73=\"],\n  [\"706\u00f76=\", \"144\u00f74=\"],\n  [\"955\u00f76=\", \"463\u00f79=\"],\n  [\"274\u00f73=\", \"325\u00f75=\"],\n  [\"631\u00f78=\", \"271\u00f77=\"],\n  [\"476\u00f72=\", \"866\u00f75=\"],\n  [\"460\u00f76=\", \"699\u00f78=\"],\n  [\"814\u00f72=\", \"853\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each \"a\u00f7b=\" exercise text with its updated value.\n# The mapping below mirrors the unified diff 1:1, in document order, and\n# every \"old\" value is unique within the document, so a simple\n# Find/Replace per pair is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-05-02 Thursday\", \"2024-05-03 Friday\"),\n  @(\"261\u00f73=\", \"250\u00f75=\"),\n  @(\"618\u00f76=\", \"977\u00f74=\"),\n  @(\"811\u00f77=\", \"931\u00f78=\"),\n  @(\"436\u00f76=\", \"446\u00f77=\"),\n  @(\"796\u00f76=\", \"380\u00f75=\"),\n  @(\"746\u00f79=\", \"818\u00f79=\"),\n  @(\"113\u00f74=\", \"715\u00f75=\"),\n  @(\"351\u00f78=\", \"332\u00f77=\"),\n  @(\"284\u00f74=\", \"496\u00f76=\"),\n  @(\"904\u00f79=\", \"492\u00f77=\"),\n  @(\"328\u00f78=\", \"293\u00f74=\"),\n  @(\"199\u00f75=\", \"960\u00f72=\"),\n  @(\"749\u00f79=\", \"749\u00f75=\"),\n  @(\"301\u00f73=\", \"192\u00f79=\"),\n  @(\"160\u00f72=\", \"133\u00f75=\"),\n  @(\"851\u00f77=\", \"181\u00f77=\"),\n  @(\"511\u00f77=\", \"571\u00f73=\"),\n  @(\"276\u00f75=\", \"165\u00f73=\"),\n  @(\"706\u00f76=\", \"144\u00f74=\"),\n  @(\"955\u00f76=\", \"463\u00f79=\"),\n  @(\"274\u00f73=\", \"325\u00f75=\"),\n  @(\"631\u00f78=\", \"271\u00f77=\"),\n  @(\"476\u00f72=\", \"866\u00f75=\"),\n  @(\"460\u00f76=\", \"699\u00f78=\"),\n  @(\"814\u00f72=\", \"853\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
